$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.55%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.84%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.04%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08386"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.88%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.807"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.65%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.010"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.74%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.465"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.904"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.19%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9242"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.15%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1281"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1976"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.64%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09546"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03843"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.50%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1062"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.07%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001300"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.08%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006104"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.56%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.427"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.75%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.66%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.825"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.52%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1362"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.09%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2508"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04399"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.00%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004378"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.75%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.74%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003989"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02874"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.30%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05527"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.71%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007957"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1435"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.01%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008996"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.70%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002070"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.07%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01167"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006933"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.36%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003463"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.56%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002278"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.20%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.14%"
